$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.893.76"
$ws.Range("E2").Value = "  +0.70%  "

$ws.Range("D3").Value = "1.551.52"
$ws.Range("E3").Value = "  +1.16%  "

$ws.Range("D5").Value = "'206.47"
$ws.Range("E5").Value = "  +0.37%  "

$ws.Range("D6").Value = "'0.486"
$ws.Range("E6").Value = "  +0.39%  "

$ws.Range("E7").Value = "  +0.49%  "

$ws.Range("D8").Value = "'21.68"
$ws.Range("E8").Value = "  +1.54%  "

$ws.Range("D9").Value = "'0.247"
$ws.Range("E9").Value = "  +1.30%  "

$ws.Range("E10").Value = "  +1.08%  "

$ws.Range("D11").Value = "'0.0859"
$ws.Range("E11").Value = "  +0.69%  "

$ws.Range("D12").Value = "1.772.87"
$ws.Range("E12").Value = "  +1.16%  "

$ws.Range("D13").Value = "1.551.22"
$ws.Range("E13").Value = "  +1.22%  "

$ws.Range("E14").Value = "  +1.46%  "

$ws.Range("D15").Value = "'0.515"
$ws.Range("E15").Value = "  +1.69%  "

$ws.Range("D16").Value = "26.891.88"
$ws.Range("E16").Value = "  +0.68%  "

$ws.Range("D17").Value = "'61.63"
$ws.Range("E17").Value = "  +0.66%  "

$ws.Range("D18").Value = "'216.73"
$ws.Range("E18").Value = "  +2.18%  "

$ws.Range("E19").Value = "  +1.04%  "

$ws.Range("E20").Value = "  +0.03%  "

$ws.Range("E21").Value = "  +0.57%  "

$ws.Range("E22").Value = "  +0.91%  "

$ws.Range("D23").Value = "'9.20"
$ws.Range("E23").Value = "  +1.32%  "

$ws.Range("E24").Value = "  +0.74%  "

$ws.Range("D25").Value = "'153.88"
$ws.Range("E25").Value = "  +1.41%  "

$ws.Range("D26").Value = "'6.59"
$ws.Range("E26").Value = "  +0.68%  "

$ws.Range("D27").Value = "'14.87"
$ws.Range("E27").Value = "  +0.28%  "

$ws.Range("E28").Value = "  +0.54%  "

$ws.Range("E29").Value = "  +0.91%  "

$ws.Range("E30").Value = "  +2.55%  "

$ws.Range("E31").Value = "  -0.45%  "

$ws.Range("E32").Value = "  -0.84%  "

$ws.Range("D33").Value = "1.426.10"
$ws.Range("E33").Value = "  +4.87%  "

$ws.Range("E34").Value = "  +2.97%  "

$ws.Range("E35").Value = "  +4.01%  "

$ws.Range("D36").Value = "'0.958"
$ws.Range("E36").Value = "  +2.22%  "

$ws.Range("E37").Value = "  +1.14%  "

$ws.Range("E38").Value = "  +0.99%  "

$ws.Range("D39").Value = "'0.521"
$ws.Range("E39").Value = "  -0.04%  "

$ws.Range("E40").Value = "  +1.27%  "

$ws.Range("E41").Value = "  +0.55%  "

$ws.Range("D42").Value = "'5.69"
$ws.Range("E42").Value = "  +0.55%  "

$ws.Range("D43").Value = "'0.985"
$ws.Range("E43").Value = "  -0.76%  "

$ws.Range("E44").Value = "  +3.83%  "

$ws.Range("D45").Value = "'63.48"
$ws.Range("E45").Value = "  +1.53%  "

$ws.Range("E46").Value = "  -0.01%  "

$ws.Range("D47").Value = "1.686.21"
$ws.Range("E47").Value = "  +1.19%  "

$ws.Range("D48").Value = "'86.15"
$ws.Range("E48").Value = "  +0.99%  "

$ws.Range("E49").Value = "  +3.39%  "

$ws.Range("E50").Value = "  +4.03%  "

$ws.Range("D51").Value = "'0.0957"
$ws.Range("E51").Value = "  +1.50%  "

